# Updated cryptos list - applies new Price (D) and Volume(1h) (E) values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("Price") values that look like plain decimal numbers must be forced
# back to text (matching the source data which stores prices as literal strings,
# e.g. thousand-dot-grouped "29.214.69"). Cells whose text already contains more
# than one "." are left alone since Excel cannot parse them as numbers anyway.
$textForceRows = @(4, 5, 6, 7, 8, 9, 10, 11, 13, 14, 15, 17, 18, 19, 20, 22, 23, 24, 25, 26, 27, 28, 29, 30, 31, 32, 34, 36, 37, 39, 40, 41, 42, 43, 44, 45, 46, 47, 48, 49, 50, 51)
foreach ($r in $textForceRows) {
    $ws.Cells.Item($r, 4).NumberFormat = "@"
}

# Column D (Price) updates
$ws.Range("D2").Value = "29.214.69"
$ws.Range("D3").Value = "1.899.17"
$ws.Range("D4").Value = "1.003"
$ws.Range("D5").Value = "326.30"
$ws.Range("D6").Value = "1.002"
$ws.Range("D7").Value = "0.4643"
$ws.Range("D8").Value = "0.3916"
$ws.Range("D9").Value = "0.07873"
$ws.Range("D10").Value = "0.9874"
$ws.Range("D11").Value = "21.93"
$ws.Range("D12").Value = "1.929.60"
$ws.Range("D13").Value = "7.070"
$ws.Range("D14").Value = "5.742"
$ws.Range("D15").Value = "0.06980"
$ws.Range("D17").Value = "1.004"
$ws.Range("D18").Value = "0.000009977"
$ws.Range("D19").Value = "17.07"
$ws.Range("D20").Value = "1.002"
$ws.Range("D21").Value = "29.234.12"
$ws.Range("D22").Value = "5.314"
$ws.Range("D23").Value = "11.07"
$ws.Range("D24").Value = "2.096"
$ws.Range("D25").Value = "156.22"
$ws.Range("D26").Value = "19.44"
$ws.Range("D27").Value = "5.970"
$ws.Range("D28").Value = "118.69"
$ws.Range("D29").Value = "1.906"
$ws.Range("D30").Value = "0.09347"
$ws.Range("D31").Value = "0.9033"
$ws.Range("D32").Value = "5.279"
$ws.Range("D34").Value = "3.218"
$ws.Range("D36").Value = "0.05775"
$ws.Range("D37").Value = "0.02084"
$ws.Range("D39").Value = "7.752"
$ws.Range("D40").Value = "0.5705"
$ws.Range("D41").Value = "0.1784"
$ws.Range("D42").Value = "9.762"
$ws.Range("D43").Value = "11.97"
$ws.Range("D44").Value = "0.5347"
$ws.Range("D45").Value = "2.186"
$ws.Range("D46").Value = "0.07043"
$ws.Range("D47").Value = "1.848"
$ws.Range("D48").Value = "2.578"
$ws.Range("D49").Value = "113.21"
$ws.Range("D50").Value = "1.050"
$ws.Range("D51").Value = "71.15"

# Column E (Volume(1h)) updates
$ws.Range("E2").Value = "  +0.14%  "
$ws.Range("E3").Value = "  -0.54%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("E5").Value = "  -0.37%  "
$ws.Range("E6").Value = "  -0.12%  "
$ws.Range("E7").Value = "  -0.35%  "
$ws.Range("E8").Value = "  -0.26%  "
$ws.Range("E9").Value = "  -1.29%  "
$ws.Range("E10").Value = "  -1.74%  "
$ws.Range("E11").Value = "  -1.67%  "
$ws.Range("E12").Value = "  +0.22%  "
$ws.Range("E13").Value = "  -0.96%  "
$ws.Range("E14").Value = "  -0.93%  "
$ws.Range("E15").Value = "  -0.05%  "
$ws.Range("E16").Value = "  -0.40%  "
$ws.Range("E17").Value = "  +0.00%  "
$ws.Range("E18").Value = "  -1.28%  "
$ws.Range("E19").Value = "  -1.06%  "
$ws.Range("E20").Value = "  -0.20%  "
$ws.Range("E21").Value = "  +0.17%  "
$ws.Range("E22").Value = "  -1.19%  "
$ws.Range("E23").Value = "  -0.14%  "
$ws.Range("E24").Value = "  +1.89%  "
$ws.Range("E25").Value = "  -0.07%  "
$ws.Range("E26").Value = "  -0.63%  "
$ws.Range("E27").Value = "  +1.99%  "
$ws.Range("E29").Value = "  -5.03%  "
$ws.Range("E30").Value = "  -0.70%  "
$ws.Range("E31").Value = "  -2.37%  "
$ws.Range("E32").Value = "  -1.65%  "
$ws.Range("E33").Value = "  -1.52%  "
$ws.Range("E34").Value = "  -1.57%  "
$ws.Range("E36").Value = "  -1.29%  "
$ws.Range("E37").Value = "  -0.74%  "
$ws.Range("E38").Value = "  -0.30%  "
$ws.Range("E39").Value = "  -3.51%  "
$ws.Range("E40").Value = "  -0.94%  "
$ws.Range("E41").Value = "  -1.50%  "
$ws.Range("E42").Value = "  -2.54%  "
$ws.Range("E43").Value = "  -0.44%  "
$ws.Range("E44").Value = "  -1.61%  "
$ws.Range("E45").Value = "  -1.91%  "
$ws.Range("E46").Value = "  -0.82%  "
$ws.Range("E47").Value = "  -2.03%  "
$ws.Range("E48").Value = "  -0.08%  "
$ws.Range("E49").Value = "  +0.70%  "
$ws.Range("E50").Value = "  -3.25%  "
$ws.Range("E51").Value = "  -0.62%  "
